$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that get a new "Yes" value in column H (Questionable), matching the
# style already used by the adjacent G column cells on the same rows.
$rows = @(38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,70)

foreach ($r in $rows) {
    $ws.Range("G$r").Copy() | Out-Null
    $ws.Range("H$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("H$r").Value = "Yes"
}
